$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Invoice sheet: add the "number part" of the Invoice Organization
#    Code to the report generator documentation - new
#    OrganizationCodeValue / OrganizationCodeDisplay rows, bracketing
#    the existing OrganizationCodeName row.
# ------------------------------------------------------------------
$invoice = $wb.Worksheets.Item("Invoice")

# Row 19 currently holds "OrganizationCodeName"; push it down one row
# and insert a fresh row above it for "OrganizationCodeValue".
$invoice.Rows.Item(19).Insert()

# "OrganizationCodeName" is now row 20; insert a fresh row below it
# (at row 21) for "OrganizationCodeDisplay".
$invoice.Rows.Item(21).Insert()

# New row 19 - OrganizationCodeValue
$invoice.Cells.Item(19, 3).Value = "OrganizationCodeValue"
$invoice.Cells.Item(19, 4).Value = "<%= invoice.OrganizationCodeValue %>"

# New row 21 - OrganizationCodeDisplay
$invoice.Cells.Item(21, 3).Value = "OrganizationCodeDisplay"
$invoice.Cells.Item(21, 4).Value = "<%= invoice.OrganizationCodeDisplay %>"

# Example values (column E) for all three related rows.
$invoice.Cells.Item(19, 5).Value = "e.g. 5900"
$invoice.Cells.Item(20, 5).Value = "e.g. Forest Resilience Division "
$invoice.Cells.Item(21, 5).Value = "e.g. 5900 - Forest Resilience Division"

# ------------------------------------------------------------------
# 2. View-state bookkeeping so the active tab moves from
#    "InvoicePaymentRequest" back to "Project", and
#    "ExpectedPerformanceMeasure" keeps its last selected cell (C74).
# ------------------------------------------------------------------
$expectedPerformanceMeasure = $wb.Worksheets.Item("ExpectedPerformanceMeasure")
$expectedPerformanceMeasure.Activate()
$expectedPerformanceMeasure.Range("C74").Select()

$project = $wb.Worksheets.Item("Project")
$project.Activate()
